$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "masterdatas"
$ws.Range("B2").Value = "testlower"
$ws.Range("C2").Value = "user"
$ws.Range("D2").Value = "password"

$ws.Range("A3").Value = "masterdatas"
$ws.Range("B3").Value = "testcamel"
$ws.Range("C3").Value = "user"
$ws.Range("D3").Value = "password"

$ws.Range("H24").Select()
